$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.179.49"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.095.37"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'229.10"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").Value = "'60.76"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.0846"
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "2.403.08"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").Value = "'14.65"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "'22.29"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").Value = "'5.49"
$ws.Range("E15").Value = "  +6.18%  "
$ws.Range("D16").Value = "'0.774"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "2.094.57"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "38.112.64"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").Value = "'70.21"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").Value = "'224.26"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("D26").Value = "'170.13"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").Value = "'9.44"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "'18.99"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").Value = "'1.37"
$ws.Range("E30").Value = "  +6.69%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  +6.91%  "
$ws.Range("D33").Value = "'4.71"
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'6.44"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'2.39"
$ws.Range("E37").Value = "  +4.83%  "
$ws.Range("D38").Value = "'3.55"
$ws.Range("E38").Value = "  +7.93%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'17.99"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "1.557.51"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").Value = "'100.25"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").Value = "'0.0219"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "'4.16"
$ws.Range("E46").Value = "  +4.54%  "
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "'7.26"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").Value = "2.290.86"
$ws.Range("E51").Value = "  +2.90%  "
